# Clean up dirs, rename main python file, fix issue with uppercase naming in risks
#
# The "Tag Matrix" sheet had two rows whose tag name was written in
# uppercase ("SSPHP-Metrics" / "SSPHP-Metrics-rust-p3sha"). They sorted to
# the top of the table (right after the header row) because of the
# capitalisation. This fixes the casing to match the rest of the tags
# ("ssphp-metrics" / "ssphp-metrics-rust-p3sha") and re-inserts the two
# rows in their correct alphabetically-sorted position (between
# "s184d01-compdefault" and "tfstatel95cd"), carrying over the same "X"
# marks each row had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "X" marks (by column letter) for the two mis-cased rows
# before we remove them.
$row2Marks = @("B", "E", "U", "AG", "AL", "AR")
$row3Marks = @("B", "C", "N", "AK", "AM")

# Remove the two uppercase rows (rows 2 and 3).
$ws.Rows("2:3").Delete()

# After the delete, the table (previously rows 4-18) now occupies rows
# 2-16:
#   2 s184d01-comp-complete-app
#   3 s184d01-comp-complete-app-worker
#   4 s184d01-comp-tfvars
#   5 s184d01-compdefault
#   6 tfstatel95cd            <- correctly-cased rows belong before this
#   7 tfstatep3sha
#   ...
# Insert two blank rows at 6 and 7 (pushing tfstatel95cd etc. back down)
# to hold the re-cased entries in alphabetical order.
$ws.Rows("6:7").Insert()

# Row 6: ssphp-metrics
$ws.Range("A6").Value = "ssphp-metrics"
foreach ($col in $row2Marks) {
    $ref = $col + "6"
    $ws.Range($ref).Value = "X"
}

# Row 7: ssphp-metrics-rust-p3sha
$ws.Range("A7").Value = "ssphp-metrics-rust-p3sha"
foreach ($col in $row3Marks) {
    $ref = $col + "7"
    $ws.Range($ref).Value = "X"
}
